# Trade #55 closed at 2026-02-17 15:43:17 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up numbers and appends the
# newly-closed trade #55 to both the "All Trades" and "MarketMaking" logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet - refresh the aggregate stats
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.37   # Current Capital
$summary.Range("B4").Value = 0.37      # Total P&L $
$summary.Range("B5").Value = 0.13      # Total P&L %
$summary.Range("B6").Value = 55        # Total Trades
$summary.Range("B7").Value = 16        # Winning Trades
$summary.Range("B9").Value = 29.09     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - refresh the MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.37     # Capital
$status.Range("D4").Value = 55         # Trades
$status.Range("E4").Value = 0.37       # P&L $
$status.Range("F4").Value = 0.37       # P&L %
$status.Range("G4").Value = 29.09      # Win Rate %

# ---------------------------------------------------------------------
# Append the new closed trade (#55) to the "All Trades" and
# "MarketMaking" logs - both sheets mirror the same data.
# ---------------------------------------------------------------------
$newRow = @{
    A = 55
    B = "2026-02-17"
    C = "15:43:10"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.66
    G = 0.7
    H = "CLOSED"
    I = 6.0606
    J = 0.04
    K = 100.37
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 56

    $ws.Range("A$row").Value = $newRow.A

    # The Date column holds plain text in this log (not a real Excel
    # date), so force a text format before assigning it to keep it as a
    # literal string instead of being auto-parsed into a date serial.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = $newRow.B
    $ws.Range("C$row").Value = $newRow.C

    $ws.Range("D$row").Value = $newRow.D
    $ws.Range("E$row").Value = $newRow.E
    $ws.Range("F$row").Value = $newRow.F
    $ws.Range("G$row").Value = $newRow.G
    $ws.Range("H$row").Value = $newRow.H
    $ws.Range("I$row").Value = $newRow.I
    $ws.Range("J$row").Value = $newRow.J
    $ws.Range("K$row").Value = $newRow.K
    $ws.Range("L$row").Value = $newRow.L
    $ws.Range("M$row").Value = $newRow.M
    $ws.Range("N$row").Value = $newRow.N
    $ws.Range("O$row").Value = $newRow.O
    $ws.Range("P$row").Value = $newRow.P
    $ws.Range("Q$row").Value = $newRow.Q
}
